$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D/E price+volume cells to remain plain text (matches source data
# which stores these as text, not numbers) while writing the new values.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.007.78'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.924.94'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('D5').Value = '325.07'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.4586'
$ws.Range('D8').Value = '0.3819'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '0.07754'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = '0.9791'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').Value = '22.57'
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = '1.945.47'
$ws.Range('E12').Value = '  +2.14%  '
$ws.Range('D13').Value = '5.711'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '6.978'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '0.06993'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').Value = '84.84'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '0.000009491'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = '16.72'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').Value = '29.016.41'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = '5.352'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('D24').Value = '2.171.94'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('D25').Value = '2.058'
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('D26').Value = '158.13'
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').Value = '19.03'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '5.630'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').Value = '117.81'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').Value = '1.841'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').Value = '0.09319'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').Value = '0.8639'
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').Value = '5.109'
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('D34').Value = '1.247'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '0.05702'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('E40').Value = '  +13.49%  '
$ws.Range('D41').Value = '7.468'
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').Value = '0.5515'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '9.363'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('D45').Value = '0.000002840'
$ws.Range('E45').Value = '  +10.86%  '
$ws.Range('D46').Value = '2.182'
$ws.Range('E46').Value = '  +4.56%  '
$ws.Range('D47').Value = '0.5185'
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.06936'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '11.19'
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('D50').Value = '111.04'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  -0.38%  '

# Restore the original (default) cell style now that the text values are set.
$textRange.Style = "Normal"
